$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13: new timesheet entry ---
# Bring over number formats / borders for the date & time columns from row 12
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E12").Copy()
$ws.Range("E13").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A13").Value = 44954
$ws.Range("B13").Value = 0.46875
$ws.Range("C13").Value = 0.51041666666666663
$ws.Range("D13").Value = "Availability and attributes can be modified. Started adding some http response codes."
$ws.Range("E13").Value = 1
$ws.Rows("13:13").RowHeight = 42

# --- Row 14: new timesheet entry ---
$ws.Range("A9:C9").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E9").Copy()
$ws.Range("E14").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A14").Value = 44971
$ws.Range("B14").Value = 0.72222222222222221
$ws.Range("C14").Value = 0.76388888888888884
$ws.Range("D14").Value = "Began work on frontend interface for adding students and tutors"
$ws.Range("E14").Value = 1
$ws.Rows("14:14").RowHeight = 28

# --- Update the selected cell to reflect where the user left off ---
$ws.Range("A15").Select()

$wb.Save()
